# Refresh the crypto price/volume table with the latest scrape values.
# (GitHub Actions scheduled update -- see commit message.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.902.68"
$ws.Range("E2").Value = "  +1.51%  "

$ws.Range("D3").Value = "1.641.01"
$ws.Range("E3").Value = "  +1.50%  "

$ws.Range("E4").Value = "  -0.16%  "

$cell = $ws.Range("D5")
$cell.Value = "'212.75"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.76%  "

$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("E7").Value = "  -0.18%  "

$cell = $ws.Range("D8")
$cell.Value = "'23.43"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +2.25%  "

$ws.Range("E9").Value = "  +2.01%  "

$ws.Range("E10").Value = "  +0.34%  "

$ws.Range("E11").Value = "  -1.90%  "

$ws.Range("D12").Value = "1.873.91"
$ws.Range("E12").Value = "  +1.51%  "

$ws.Range("D13").Value = "1.644.56"
$ws.Range("E13").Value = "  +1.64%  "

$ws.Range("E14").Value = "  +1.25%  "

$ws.Range("E15").Value = "  +2.64%  "

$cell = $ws.Range("D16")
$cell.Value = "'65.64"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +2.06%  "

$ws.Range("D17").Value = "27.891.46"
$ws.Range("E17").Value = "  +1.44%  "

$cell = $ws.Range("D18")
$cell.Value = "'231.54"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.20%  "

$cell = $ws.Range("D19")
$cell.Value = "'7.68"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.72%  "

$ws.Range("E20").Value = "  +0.48%  "

$ws.Range("E21").Value = "  -0.14%  "

$cell = $ws.Range("D22")
$cell.Value = "'10.74"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +8.77%  "

$cell = $ws.Range("D23")
$cell.Value = "'4.38"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +2.28%  "

$ws.Range("E24").Value = "  +3.70%  "

$cell = $ws.Range("D25")
$cell.Value = "'151.67"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.79%  "

$ws.Range("E26").Value = "  +1.06%  "

$cell = $ws.Range("D27")
$cell.Value = "'0.111"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.21%  "

$cell = $ws.Range("D28")
$cell.Value = "'15.69"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.97%  "

$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("E30").Value = "  +0.43%  "

$ws.Range("E31").Value = "  +0.35%  "

$ws.Range("E32").Value = "  +0.99%  "

$ws.Range("D33").Value = "1.456.49"
$ws.Range("E33").Value = "  +0.69%  "

$ws.Range("E34").Value = "  +1.38%  "

$cell = $ws.Range("D35")
$cell.Value = "'1.55"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +1.32%  "

$ws.Range("E36").Value = "  -0.53%  "

$cell = $ws.Range("D37")
$cell.Value = "'0.889"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +3.33%  "

$cell = $ws.Range("D38")
$cell.Value = "'0.563"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +0.25%  "

$ws.Range("E39").Value = "  +0.95%  "

$cell = $ws.Range("D40")
$cell.Value = "'0.916"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.96%  "

$cell = $ws.Range("D41")
$cell.Value = "'69.27"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.49%  "

$ws.Range("E43").Value = "  +0.49%  "

$cell = $ws.Range("D44")
$cell.Value = "'2.48"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("E45").Value = "  +0.78%  "

$ws.Range("E46").Value = "  -0.27%  "

$ws.Range("E47").Value = "  +6.09%  "

$ws.Range("D48").Value = "1.782.95"
$ws.Range("E48").Value = "  +1.29%  "

$cell = $ws.Range("D49")
$cell.Value = "'88.42"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +2.74%  "

$ws.Range("E50").Value = "  +0.26%  "

$ws.Range("E51").Value = "  +2.24%  "
